$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H53").Value = 1745.579
$ws.Range("I53").Value = 192
$ws.Range("J53").Value = 2651.8333
$ws.Range("K53").Value = 192
$ws.Range("L53").Value = 2651.8333
$ws.Range("M53").Value = 445
$ws.Range("N53").Value = -3925.8333
$ws.Range("H86").Value = 21048.4
$ws.Range("J86").Value = 50521
$ws.Range("L86").Value = 50521
$ws.Range("N86").Value = -52767
$ws.Range("H89").Value = 21048.4
$ws.Range("J89").Value = 50521
$ws.Range("L89").Value = 252605
$ws.Range("N89").Value = -263837
$ws.Range("H113").Value = 40004720
$ws.Range("I113").Value = 76927000
$ws.Range("J113").Value = 5583.3335
$ws.Range("K113").Value = 76927000
$ws.Range("L113").Value = 5583.3335
$ws.Range("M113").Value = -76923746
$ws.Range("N113").Value = -12091.3335
$ws.Range("H137").Value = 24320.418
$ws.Range("I137").Value = 805.93335
$ws.Range("J137").Value = 78584.62
$ws.Range("K137").Value = 2417.80005
$ws.Range("L137").Value = 235753.86
$ws.Range("M137").Value = 132.1999500000002
$ws.Range("N137").Value = -240853.86
$ws.Range("H141").Value = 1116.6595
$ws.Range("I141").Value = 854.4474
$ws.Range("J141").Value = 2223.7778
$ws.Range("K141").Value = 2563.3422
$ws.Range("L141").Value = 6671.3334
$ws.Range("M141").Value = 2616.6578
$ws.Range("N141").Value = -17031.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15382.161
$ws.Range("I32").Value = 16883.04
$ws.Range("J32").Value = 5012.4546
$ws.Range("K32").Value = 16883.04
$ws.Range("L32").Value = 5012.4546
$ws.Range("M32").Value = -16596.04
$ws.Range("N32").Value = -5586.4546
$ws.Range("H74").Value = 18183056
$ws.Range("I74").Value = 20409484
$ws.Range("J74").Value = 564.6667
$ws.Range("K74").Value = 20409484
$ws.Range("L74").Value = 564.6667
$ws.Range("M74").Value = -20408610
$ws.Range("N74").Value = -2312.6667
$ws.Range("H77").Value = 18183056
$ws.Range("I77").Value = 20409484
$ws.Range("J77").Value = 564.6667
$ws.Range("K77").Value = 102047420
$ws.Range("L77").Value = 2823.3335
$ws.Range("M77").Value = -102043052
$ws.Range("N77").Value = -11559.3335
$ws.Range("H132").Value = 42291.594
$ws.Range("I132").Value = 1811.8572
$ws.Range("J132").Value = 168228.56
$ws.Range("K132").Value = 5435.571599999999
$ws.Range("L132").Value = 504685.68
$ws.Range("M132").Value = -2905.571599999999
$ws.Range("N132").Value = -509745.68

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1267.0667
$ws.Range("I94").Value = 807.38464
$ws.Range("K94").Value = 807.38464
$ws.Range("M94").Value = -356.38464
$ws.Range("H134").Value = 40273.395
$ws.Range("I134").Value = 46360.625
$ws.Range("K134").Value = 139081.875
$ws.Range("M134").Value = -136546.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8256.25
$ws.Range("I31").Value = 8344.768
$ws.Range("J31").Value = 7833.3335
$ws.Range("K31").Value = 8344.768
$ws.Range("L31").Value = 7833.3335
$ws.Range("M31").Value = -8049.768
$ws.Range("N31").Value = -8423.333500000001
$ws.Range("H34").Value = 8256.25
$ws.Range("I34").Value = 8344.768
$ws.Range("J34").Value = 7833.3335
$ws.Range("K34").Value = 8344.768
$ws.Range("L34").Value = 7833.3335
$ws.Range("M34").Value = -8142.768
$ws.Range("N34").Value = -8237.333500000001
$ws.Range("H58").Value = 10033.709
$ws.Range("I58").Value = 921.09753
$ws.Range("J58").Value = 36720.645
$ws.Range("K58").Value = 921.09753
$ws.Range("L58").Value = 36720.645
$ws.Range("M58").Value = -718.09753
$ws.Range("N58").Value = -37126.645
$ws.Range("H132").Value = 10163.966
$ws.Range("I132").Value = 10399.667
$ws.Range("K132").Value = 31199.001
$ws.Range("M132").Value = -28669.001
$ws.Range("H134").Value = 743.2368
$ws.Range("I134").Value = 720.6389
$ws.Range("K134").Value = 2161.9167
$ws.Range("M134").Value = 373.0832999999998
$ws.Range("H136").Value = 10033.709
$ws.Range("I136").Value = 921.09753
$ws.Range("J136").Value = 36720.645
$ws.Range("K136").Value = 2763.29259
$ws.Range("L136").Value = 110161.935
$ws.Range("M136").Value = -213.29259
$ws.Range("N136").Value = -115261.935

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5100
$ws.Range("J63").Value = 9500
$ws.Range("L63").Value = 28500
$ws.Range("N63").Value = -29998
$ws.Range("H66").Value = 5100
$ws.Range("J66").Value = 9500
$ws.Range("L66").Value = 85500
$ws.Range("N66").Value = -92988
$ws.Range("H92").Value = 10417019
$ws.Range("I92").Value = 25000270
$ws.Range("J92").Value = 410.7143
$ws.Range("K92").Value = 75000810
$ws.Range("L92").Value = 1232.1429
$ws.Range("M92").Value = -74999562
$ws.Range("N92").Value = -3728.1429
$ws.Range("H131").Value = 704.29
$ws.Range("J131").Value = 716.76044
$ws.Range("L131").Value = 2150.28132
$ws.Range("N131").Value = -12230.28132
$ws.Range("H140").Value = 2599.5715
$ws.Range("I140").Value = 2056.8572
$ws.Range("K140").Value = 6170.571599999999
$ws.Range("M140").Value = -990.5715999999993

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3520.75
$ws.Range("I80").Value = 2782.6667
$ws.Range("K80").Value = 2782.6667
$ws.Range("M80").Value = -1784.6667
$ws.Range("H83").Value = 3520.75
$ws.Range("I83").Value = 2782.6667
$ws.Range("K83").Value = 13913.3335
$ws.Range("M83").Value = -8921.333500000001
$ws.Range("H122").Value = 95239176
$ws.Range("I122").Value = 47620050
$ws.Range("J122").Value = 142858300
$ws.Range("K122").Value = 142860150
$ws.Range("L122").Value = 428574900
$ws.Range("M122").Value = -142857700
$ws.Range("N122").Value = -428579800

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2084.9285
$ws.Range("I100").Value = 1580.6
$ws.Range("J100").Value = 2365.111
$ws.Range("K100").Value = 1580.6
$ws.Range("L100").Value = 2365.111
$ws.Range("M100").Value = -1039.6
$ws.Range("N100").Value = -3447.111
$ws.Range("H140").Value = 50134.25
$ws.Range("J140").Value = 50134.25
$ws.Range("L140").Value = 50134.25
$ws.Range("N140").Value = -60494.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 587.3043
$ws.Range("I132").Value = 535.4
$ws.Range("K132").Value = 1606.2
$ws.Range("M132").Value = 923.8000000000002

Write-Output "Applied all Typhon_Profits updates"